$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the two new columns -------------------------------------
# "Sub Category" goes right after "Category" (old column B).
$ws.Columns("B").Insert()
# "Supplier" goes right after "UOM". Before this second insert, UOM is
# still column K (A,B(new),C..K=UOM), so inserting at L puts the new
# column right after it.
$ws.Columns("L").Insert()

# --- 2. Insert a new row for the extra equipment entry ------------------
# This pushes the existing data row (row 3) down to row 4.
$ws.Rows("3").Insert()

# Give the freshly inserted row the same formatting as the data row below
# it (row 4) before filling in values.
$ws.Range("A4:X4").Copy($ws.Range("A3:X3"))

Write-Output "structure done"
